$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New date for this auto-update run (keep as text, not an auto-converted date serial)
$newDate = "2025-12-03"
$dateCells = $ws.Range("A2:A5")
$dateCells.NumberFormat = "@"
$dateCells.Value = $newDate
$dateCells.Style = "Normal"

# Row 2: HD HYUNDAI MIPO / 010620.KS
$ws.Cells.Item(2, 2).Value = "HD HYUNDAI MIPO"
$ws.Cells.Item(2, 3).Value = "010620.KS"
$ws.Cells.Item(2, 4).Value = 223000
$ws.Cells.Item(2, 5).Value = 60
$ws.Cells.Item(2, 6).Value = 0.68
$ws.Cells.Item(2, 7).Value = 50
$ws.Cells.Item(2, 8).Value = 53
$ws.Cells.Item(2, 9).Value = 56
$ws.Cells.Item(2, 10).Value = 66
$ws.Cells.Item(2, 11).Value = 57
$ws.Cells.Item(2, 14).Value = 65.32892478746797

# Row 3: HDKSOE / 009540.KS
$ws.Cells.Item(3, 2).Value = "HDKSOE"
$ws.Cells.Item(3, 3).Value = "009540.KS"
$ws.Cells.Item(3, 4).Value = 413000
$ws.Cells.Item(3, 5).Value = 45.5
$ws.Cells.Item(3, 6).Value = -1.67
$ws.Cells.Item(3, 7).Value = 20
$ws.Cells.Item(3, 8).Value = 63
$ws.Cells.Item(3, 9).Value = 70
$ws.Cells.Item(3, 10).Value = 76
$ws.Cells.Item(3, 11).Value = 53.6
$ws.Cells.Item(3, 14).Value = 65.32892478746797

# Row 4: Hanwha Ocean / 042660.KS
$ws.Cells.Item(4, 2).Value = "Hanwha Ocean"
$ws.Cells.Item(4, 3).Value = "042660.KS"
$ws.Cells.Item(4, 4).Value = 106100
$ws.Cells.Item(4, 5).Value = 24.5
$ws.Cells.Item(4, 6).Value = -6.11
$ws.Cells.Item(4, 7).Value = 10
$ws.Cells.Item(4, 8).Value = 63
$ws.Cells.Item(4, 9).Value = 76
$ws.Cells.Item(4, 10).Value = 83
$ws.Cells.Item(4, 11).Value = 53
$ws.Cells.Item(4, 14).Value = 65.32892478746797

# Row 5: SamsungHvyInd / 010140.KS
$ws.Cells.Item(5, 2).Value = "SamsungHvyInd"
$ws.Cells.Item(5, 3).Value = "010140.KS"
$ws.Cells.Item(5, 4).Value = 24500
$ws.Cells.Item(5, 5).Value = 36
$ws.Cells.Item(5, 6).Value = -1.21
$ws.Cells.Item(5, 7).Value = 10
$ws.Cells.Item(5, 8).Value = 66
$ws.Cells.Item(5, 9).Value = 76
$ws.Cells.Item(5, 10).Value = 90
$ws.Cells.Item(5, 11).Value = 53
$ws.Cells.Item(5, 14).Value = 65.32892478746797

# MACRO_SIGNAL text update (column O) for all data rows
$ws.Range("O2:O5").Value = "🟢 상승 우위 (다소 완화)"
